# Apply edits described by the diff:
# 1. Rename the worksheet from "Scanner" to "Psychiatry"
# 2. Shift the "Log Time" values (column D, rows 2-38) from 12:xx:xx to 11:xx:xx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet
$ws.Name = "Psychiatry"

# 2. Update the Log Time column values (stored as text strings "HH:MM:SS")
for ($row = 2; $row -le 38; $row++) {
    $cell = $ws.Cells.Item($row, 4)  # Column D
    $current = [string]$cell.Value2
    if ($current.StartsWith("12:")) {
        $cell.Value2 = "11:" + $current.Substring(3)
    }
}
